$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "OpenAI Embeddings" cosine-similarity row (row 5) with the
# refreshed values.
$ws.Range("B5").Value = 77.892614657367901
$ws.Range("C5").Value = 68.593206345652106
$ws.Range("D5").Value = 87.923206775194302
$ws.Range("E5").Value = 68.241956407044299
$ws.Range("F5").Value = 73.752703185625904
$ws.Range("G5").Value = 61.259439894891997

# Re-enter the "Avg" formula row (row 6) and the Lang-Avg formula row
# (row 19) as filled/shared formulas across the row.
$ws.Range("B6:H6").Formula = "=AVERAGE(B3:B5)"
$ws.Range("B19:H19").Formula = "=AVERAGE(B10:B18)"

# Move the active selection.
$ws.Range("J10").Select() | Out-Null
